$wb = $excel.ActiveWorkbook

# --- Sheet "Persone-Indisp" ---
$ws1 = $wb.Worksheets.Item("Persone-Indisp")

# A2: "BAI" -> "AAA"
$ws1.Range("A2").Value = "AAA"
# Remove old "G"/"GN" helper cells that used to sit next to A2/A3
$ws1.Range("B2").Clear()
$ws1.Range("C2").Clear()
$ws1.Range("C3").Clear()
$ws1.Range("L3").Clear()
# Remove the "N" row that used to span B5:AF5
$ws1.Range("B5:AF5").Clear()
# A12: "AIN" -> "ZZZ"
$ws1.Range("A12").Value = "ZZZ"

# New formatted (underlined) empty cell at G11
$ws1.Range("G11").Font.Underline = 2
$ws1.Range("G11").HorizontalAlignment = -4108
$ws1.Range("G11").VerticalAlignment = -4108

# Selection left on L29 before switching away from this sheet
$ws1.Range("L29").Select()

# --- Sheet "Turni Fissi" ---
$ws2 = $wb.Worksheets.Item("Turni Fissi")

# Clear out the "punteggio" column (B) for every data row; underline the
# (now blank) cells as part of the new score-formula formatting.
$scoreRange = $ws2.Range("B3:B32")
$scoreRange.Font.Underline = 2
$scoreRange.HorizontalAlignment = -4108
$scoreRange.VerticalAlignment = -4108
$scoreRange.ClearContents()

# Make "Turni Fissi" the active sheet/tab with B3:B32 selected
$ws2.Activate()
$ws2.Range("B3:B32").Select()
